$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 338, pushing existing rows 338-356 down to 340-358
$ws.Range("A338:A339").EntireRow.Insert()

# Row 338: new weekly record (Primera)
$ws.Cells.Item(338, 1).Value = 4
$ws.Cells.Item(338, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(338, 3).Value = "Los Lagos"
$ws.Cells.Item(338, 4).Value = 44746
$ws.Cells.Item(338, 5).Value = 10
$ws.Cells.Item(338, 6).Value = 100112023
$ws.Cells.Item(338, 7).Value = "Brócoli"
$ws.Cells.Item(338, 8).Value = "Sin especificar"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 500
$ws.Cells.Item(338, 11).Value = 1600
$ws.Cells.Item(338, 12).Value = 1600
$ws.Cells.Item(338, 13).Value = 1600
$ws.Cells.Item(338, 14).Value = "$/unidad"
$ws.Cells.Item(338, 15).Value = "Región del Maule"
$ws.Cells.Item(338, 16).Value = 1600
$ws.Cells.Item(338, 17).Value = 1
$ws.Cells.Item(338, 18).Value = "Hortaliza"

# Row 339: new weekly record (Segunda)
$ws.Cells.Item(339, 1).Value = 4
$ws.Cells.Item(339, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(339, 3).Value = "Los Lagos"
$ws.Cells.Item(339, 4).Value = 44746
$ws.Cells.Item(339, 5).Value = 10
$ws.Cells.Item(339, 6).Value = 100112023
$ws.Cells.Item(339, 7).Value = "Brócoli"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Segunda"
$ws.Cells.Item(339, 10).Value = 250
$ws.Cells.Item(339, 11).Value = 1200
$ws.Cells.Item(339, 12).Value = 1200
$ws.Cells.Item(339, 13).Value = 1200
$ws.Cells.Item(339, 14).Value = "$/unidad"
$ws.Cells.Item(339, 15).Value = "Región del Maule"
$ws.Cells.Item(339, 16).Value = 1200
$ws.Cells.Item(339, 17).Value = 1
$ws.Cells.Item(339, 18).Value = "Hortaliza"
